# Update cryptos list snapshot values (prices + 1h volume %) and re-rank
# Solana / Cardano / Dogecoin rows, per the scheduled GitHub Actions refresh.
#
# NOTE: several "Price" column values look like plain numbers (e.g. "215.70",
# "0.0620") but must stay literal TEXT (matching the existing inlineStr cells)
# so formatting like trailing/leading zeros survives. Assigning such a string
# straight to .Value lets Excel auto-coerce it to a number, so for those
# cells we briefly force a text number-format, write the value, then restore
# the original "General" / default style so nothing else about the cell
# changes.

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "27.012.26"
$ws.Range("E2").Value = "  +0.22%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "1.683.31"
$ws.Range("E3").Value = "  +0.48%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.01%  "

# --- Row 5 (BNB) ---
Set-TextValue $ws.Range("D5") "215.70"
$ws.Range("E5").Value = "  -0.15%  "

# --- Row 6 (XRP) ---
Set-TextValue $ws.Range("D6") "0.517"
$ws.Range("E6").Value = "  -2.57%  "

# --- Row 7 (USDC) ---
$ws.Range("E7").Value = "  +0.02%  "

# --- Row 8 : was Cardano, now Solana ---
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D8") "21.46"
$ws.Range("E8").Value = "  +5.20%  "

# --- Row 9 : was Dogecoin, now Cardano ---
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D9") "0.251"
$ws.Range("E9").Value = "  -1.54%  "

# --- Row 10 : was Solana, now Dogecoin ---
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D10") "0.0620"
$ws.Range("E10").Value = "  -0.27%  "

# --- Row 11 (TRON) ---
$ws.Range("E11").Value = "  -0.40%  "

# --- Row 12 (WrappedliquidstakedEther2.0) ---
$ws.Range("D12").Value = "1.920.13"
$ws.Range("E12").Value = "  +0.36%  "

# --- Row 13 (WrappedEther) ---
$ws.Range("D13").Value = "1.663.38"
$ws.Range("E13").Value = "  -2.02%  "

# --- Row 14 (Polkadot) ---
$ws.Range("E14").Value = "  +0.20%  "

# --- Row 15 (Polygon) ---
Set-TextValue $ws.Range("D15") "0.534"
$ws.Range("E15").Value = "  +2.10%  "

# --- Row 16 (Litecoin) ---
Set-TextValue $ws.Range("D16") "65.99"

# --- Row 17 (WrappedBTC) ---
$ws.Range("D17").Value = "27.033.71"

# --- Row 18 (Chainlink) ---
Set-TextValue $ws.Range("D18") "8.15"
$ws.Range("E18").Value = "  +3.77%  "

# --- Row 19 (BitcoinCash) ---
Set-TextValue $ws.Range("D19") "236.68"
$ws.Range("E19").Value = "  +1.46%  "

# --- Row 20 (ShibaInu) ---
$ws.Range("E20").Value = "  -0.08%  "

# --- Row 22 (Uniswap) ---
Set-TextValue $ws.Range("D22") "4.45"
$ws.Range("E22").Value = "  -0.42%  "

# --- Row 23 (Avalanche) ---
Set-TextValue $ws.Range("D23") "9.24"
$ws.Range("E23").Value = "  +0.14%  "

# --- Row 24 (Toncoin) ---
Set-TextValue $ws.Range("D24") "2.13"
$ws.Range("E24").Value = "  -4.13%  "

# --- Row 25 (Monero) ---
Set-TextValue $ws.Range("D25") "146.70"
$ws.Range("E25").Value = "  +0.63%  "

# --- Row 26 (Cosmos) ---
Set-TextValue $ws.Range("D26") "7.24"
$ws.Range("E26").Value = "  +1.08%  "

# --- Row 27 (EthereumClassic) ---
Set-TextValue $ws.Range("D27") "16.07"
$ws.Range("E27").Value = "  +0.62%  "

# --- Row 28 (Stellar) ---
$ws.Range("E28").Value = "  -2.82%  "

# --- Row 29 (BinanceUSD) ---
$ws.Range("E29").Value = "  +0.00%  "

# --- Row 30 (Hedera) ---
Set-TextValue $ws.Range("D30") "0.0501"
$ws.Range("E30").Value = "  +0.57%  "

# --- Row 31 (PancakeSwap) ---
$ws.Range("E31").Value = "  -0.48%  "

# --- Row 32 (Filecoin) ---
$ws.Range("E32").Value = "  +0.13%  "

# --- Row 33 (Maker) ---
$ws.Range("D33").Value = "1.502.20"
$ws.Range("E33").Value = "  +2.74%  "

# --- Row 34 (InternetComputer(DFINITY)) ---
$ws.Range("E34").Value = "  +0.24%  "

# --- Row 35 (LidoDAOToken) ---
$ws.Range("E35").Value = "  +4.26%  "

# --- Row 37 (ImmutableX) ---
Set-TextValue $ws.Range("D37") "0.587"
$ws.Range("E37").Value = "  +2.96%  "

# --- Row 38 (ARBITRUM) ---
Set-TextValue $ws.Range("D38") "0.917"
$ws.Range("E38").Value = "  +0.95%  "

# --- Row 39 (VeChain) ---
$ws.Range("E39").Value = "  +3.31%  "

# --- Row 40 (WEMIXToken) ---
$ws.Range("E40").Value = "  +7.38%  "

# --- Row 41 (FraxShare) ---
$ws.Range("E41").Value = "  -4.69%  "

# --- Row 42 (PaxDollar) ---
$ws.Range("E42").Value = "  +0.07%  "

# --- Row 43 (Aave) ---
Set-TextValue $ws.Range("D43") "67.76"
$ws.Range("E43").Value = "  +2.73%  "

# --- Row 44 (MXToken) ---
$ws.Range("E44").Value = "  -1.03%  "

# --- Row 45 (RocketPoolETH) ---
$ws.Range("D45").Value = "1.826.82"
$ws.Range("E45").Value = "  +0.42%  "

# --- Row 46 (TrustWalletToken) ---
Set-TextValue $ws.Range("D46") "0.780"
$ws.Range("E46").Value = "  -0.19%  "

# --- Row 47 (Quant) ---
Set-TextValue $ws.Range("D47") "90.35"
$ws.Range("E47").Value = "  -0.39%  "

# --- Row 48 (RenderToken) ---
$ws.Range("E48").Value = "  -0.51%  "

# --- Row 49 (Algorand) ---
$ws.Range("E49").Value = "  +3.65%  "

# --- Row 50 (EnergySwap) ---
Set-TextValue $ws.Range("D50") "7.85"
$ws.Range("E50").Value = "  +3.01%  "

# --- Row 51 (Cronos) ---
$ws.Range("E51").Value = "  +0.10%  "
